$d = $word.ActiveDocument
$BR = [char]11

# ---------------------------------------------------------------------------
# This document has 7 "Heading 2" sections (Objetivos, Docente(s)
# Responsavel(eis), Programa resumido, Programa, Avaliacao, Bibliografia,
# Requisitos). The edit rotates the body content of the first six sections
# (Requisitos is untouched). Paragraph count/order does not change except
# that paragraph 8 (the body under "Docente(s) Responsavel(eis)") gains a
# second run; that is handled with a temporary paragraph split + merge so
# the two text blocks stay in separate <w:r> runs (matching the target
# canonical XML), and the paragraph count is restored to 18 afterwards.
# ---------------------------------------------------------------------------

# --- Paragraph 6 (body of "Objetivos") ------------------------------------
# old: "Passar aos alunos ..."
# new: the old "Programa resumido" numbered list (no space after the dot)
$p6 = $d.Paragraphs.Item(6)
$old6 = "Passar aos alunos os conhecimentos básicos da estrutura de um vidro, a influência da composição nas propriedades e os processos empregados na produção de vidros"
$new6 = "1.Introdução, Quadro da indústria brasileira de vidros" + $BR + `
        "2.Composição dos vidros" + $BR + `
        "3.Materiais Primas" + $BR + `
        "4.Mecanismo de fusão e formação do vidro" + $BR + `
        "5.Viscosidade - Definição, relação com a composição, métodos experimentais de medição, cálculo a partir da composição" + $BR + `
        "6.Propriedades óticas " + $BR + `
        "7.Propriedades mecânicas " + $BR + `
        "8.Propriedades químicas " + $BR + `
        "9.Processamento - Vidro plano, vidro oco, vidros especiais, vidro temperado, esmalte" + $BR + `
        "10.Aula prática - Fundir um vidro, produzir um vidro colorido e esmaltar um metal"
$p6.Range.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null

# --- Paragraph 8 (body of "Docente(s) Responsável(eis)") ------------------
# old: two runs "5983729 - Fernando Vernilli Junior" <br/> | "1922320 - Sebastiao Ribeiro"
# new: two runs "Passar aos alunos ..." <br/> | "1. Introdução ... 10. Aula prática ..."
$p8 = $d.Paragraphs.Item(8)
$p8Start = $p8.Range.Start
$p8End = $p8.Range.End
$clearRng = $d.Range($p8Start, $p8End - 1)
$clearRng.Text = ""

$seg1 = "Passar aos alunos os conhecimentos básicos da estrutura de um vidro, a influência da composição nas propriedades e os processos empregados na produção de vidros" + $BR
$p8 = $d.Paragraphs.Item(8)
$insPt = $d.Range($p8.Range.Start, $p8.Range.Start)
$insPt.InsertAfter($seg1)

$p8 = $d.Paragraphs.Item(8)
$endPt = $d.Range($p8.Range.End - 1, $p8.Range.End - 1)
$endPt.InsertParagraphAfter()

$seg2 = "1. Introdução, Quadro da indústria brasileira de vidros" + $BR + `
        "2. Composição dos vidros" + $BR + `
        "3. Materiais Primas" + $BR + `
        "4. Mecanismo de fusão e formação do vidro" + $BR + `
        "5. Viscosidade – Definição, relação com a composição, métodos experimentais de medição, cálculo a partir da composição" + $BR + `
        "6. Propriedades óticas " + $BR + `
        "7. Propriedades mecânicas " + $BR + `
        "8. Propriedades químicas " + $BR + `
        "9. Processamento – Vidro plano, vidro oco, vidros especiais, vidro temperado, esmalte" + $BR + `
        "10. Aula prática - Fundir um vidro, produzir um vidro colorido e esmaltar um metal"
$p9 = $d.Paragraphs.Item(9)
$ins2 = $d.Range($p9.Range.Start, $p9.Range.Start)
$ins2.InsertAfter($seg2)

$p8 = $d.Paragraphs.Item(8)
$markStart = $p8.Range.End - 1
$markRng = $d.Range($markStart, $markStart + 1)
$markRng.Delete() | Out-Null

# --- Paragraph 10 (body of "Programa resumido") ----------------------------
# old: numbered list (no space after dot)
# new: "Serão realizadas duas provas ..." (old "Método" value of Avaliação)
$p10 = $d.Paragraphs.Item(10)
$old10 = "1.Introdução, Quadro da indústria brasileira de vidros" + $BR + `
        "2.Composição dos vidros" + $BR + `
        "3.Materiais Primas" + $BR + `
        "4.Mecanismo de fusão e formação do vidro" + $BR + `
        "5.Viscosidade - Definição, relação com a composição, métodos experimentais de medição, cálculo a partir da composição" + $BR + `
        "6.Propriedades óticas " + $BR + `
        "7.Propriedades mecânicas " + $BR + `
        "8.Propriedades químicas " + $BR + `
        "9.Processamento - Vidro plano, vidro oco, vidros especiais, vidro temperado, esmalte" + $BR + `
        "10.Aula prática - Fundir um vidro, produzir um vidro colorido e esmaltar um metal"
$new10 = "Serão realizadas duas provas escritas (P1 e P2), apresentações orais de trabalhos (T) e listas de exercícios (E)"
$p10.Range.Find.Execute($old10, $true, $false, $false, $false, $false, $true, 1, $false, $new10, 2) | Out-Null

# --- Paragraph 12 (body of "Programa") -------------------------------------
# old: numbered list (with space after dot)
# new: "A nota final será calculada ..." (old "Critério" value of Avaliação)
$p12 = $d.Paragraphs.Item(12)
$old12 = "1. Introdução, Quadro da indústria brasileira de vidros" + $BR + `
        "2. Composição dos vidros" + $BR + `
        "3. Materiais Primas" + $BR + `
        "4. Mecanismo de fusão e formação do vidro" + $BR + `
        "5. Viscosidade – Definição, relação com a composição, métodos experimentais de medição, cálculo a partir da composição" + $BR + `
        "6. Propriedades óticas " + $BR + `
        "7. Propriedades mecânicas " + $BR + `
        "8. Propriedades químicas " + $BR + `
        "9. Processamento – Vidro plano, vidro oco, vidros especiais, vidro temperado, esmalte" + $BR + `
        "10. Aula prática - Fundir um vidro, produzir um vidro colorido e esmaltar um metal"
$new12 = "A nota final será calculada utilizando a equação: {[(P1 + P2 + T)/3] x 0,9} + E x 0,1"
$p12.Range.Find.Execute($old12, $true, $false, $false, $false, $false, $true, 1, $false, $new12, 2) | Out-Null

# --- Paragraph 14 (body of "Avaliação") ------------------------------------
# Three labelled values; only the value runs change (label runs untouched).
# NOTE: the new "Método" value is textually identical to the old "Norma de
# recuperação" value, so the "Norma de recuperação" replacement MUST happen
# first -- otherwise the later Método replacement's output would be found
# (and wrongly re-replaced) by the Norma step.

# Norma de recuperação value: old "Para a recuperação ..." -> new "5983729 - Fernando Vernilli Junior"
$p14 = $d.Paragraphs.Item(14)
$oldNorma = "Para a recuperação será realizada uma prova  (PR) abrangendo toda a matéria no semestre, valendo de 0 (zero) a 10 (10). Média Final: (MP + PR)/2. Média Final igual ou superior a 5: aprovado. Média Final inferior a 5: reprovado"
$newNorma = "5983729 - Fernando Vernilli Junior"
$p14.Range.Find.Execute($oldNorma, $true, $false, $false, $false, $false, $true, 1, $false, $newNorma, 2) | Out-Null

# Critério value: old "A nota final ..." -> new bibliography list (4 lines)
$p14 = $d.Paragraphs.Item(14)
$oldCriterio = "A nota final será calculada utilizando a equação: {[(P1 + P2 + T)/3] x 0,9} + E x 0,1"
$newCriterio = "1.)Associação Brasileira da Industria de Vidros, www.abividro.br" + $BR + `
               "2.)H. Scholze, Glas, Springer-Verlag, 1988" + $BR + `
               "3.)R. H. Doremus, Glass Science, New York, John Wiley, 1994" + $BR + `
               "4.)H. G. Pfaender, Schott Guide to Glass, London, Chapman & Hall, 1996"
$p14.Range.Find.Execute($oldCriterio, $true, $false, $false, $false, $false, $true, 1, $false, $newCriterio, 2) | Out-Null

# Método value: old "Serão realizadas ..." -> new "Para a recuperação ..."
$p14 = $d.Paragraphs.Item(14)
$oldMetodo = "Serão realizadas duas provas escritas (P1 e P2), apresentações orais de trabalhos (T) e listas de exercícios (E)"
$newMetodo = "Para a recuperação será realizada uma prova  (PR) abrangendo toda a matéria no semestre, valendo de 0 (zero) a 10 (10). Média Final: (MP + PR)/2. Média Final igual ou superior a 5: aprovado. Média Final inferior a 5: reprovado"
$p14.Range.Find.Execute($oldMetodo, $true, $false, $false, $false, $false, $true, 1, $false, $newMetodo, 2) | Out-Null

# --- Paragraph 16 (body of "Bibliografia") ----------------------------------
# old: 4-line bibliography list
# new: "1922320 - Sebastiao Ribeiro"
$p16 = $d.Paragraphs.Item(16)
$old16 = "1.)Associação Brasileira da Industria de Vidros, www.abividro.br" + $BR + `
         "2.)H. Scholze, Glas, Springer-Verlag, 1988" + $BR + `
         "3.)R. H. Doremus, Glass Science, New York, John Wiley, 1994" + $BR + `
         "4.)H. G. Pfaender, Schott Guide to Glass, London, Chapman & Hall, 1996"
$new16 = "1922320 - Sebastiao Ribeiro"
$p16.Range.Find.Execute($old16, $true, $false, $false, $false, $false, $true, 1, $false, $new16, 2) | Out-Null

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
